# Eva QI_timetable.xlsx - apply timetable generation logic update:
#  - exclude activities after 17:00 for Days 1-5 (remove the 17:15 "Free Time"
#    row entries and shrink the Day blocks that used to run through to 19:00)
#  - rename two private-lesson entries to "First Last" order
#  - turn the two "Free Time" slots at 11:00/14:15 on Day 5 into the new
#    "Master class with Ivy & Stephane" joint session
#  - update the affected column widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column width changes (B, D, F, J)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth  = 34.17   # B: 36 -> 35
$ws.Columns.Item(4).ColumnWidth  = 31.17   # D: 33 -> 32
$ws.Columns.Item(6).ColumnWidth  = 34.17   # F: 36 -> 35
$ws.Columns.Item(10).ColumnWidth = 33.17   # J: 14 -> 34

# ---------------------------------------------------------------------------
# 2. Cell text updates
# ---------------------------------------------------------------------------
$ws.Range("B7").Value  = "Private lesson with Stephane RETY"
$ws.Range("J7").Value  = "Master class with Ivy & Stephane"
$ws.Range("F11").Value = "Private lesson with Stephane RETY"
$ws.Range("D20").Value = "Private lesson with Ivy CHUANG"
$ws.Range("J20").Value = "Master class with Ivy & Stephane"

# ---------------------------------------------------------------------------
# 3. Remove the 17:15 "Free Time" row entries (row 32) entirely - clear
#    formatting before clearing the value so the cell element itself drops
#    out of the saved sheet instead of lingering as an empty styled cell.
# ---------------------------------------------------------------------------
foreach ($col in @("B", "D", "F", "H", "J")) {
    $cell = $ws.Range($col + "32")
    $cell.ClearFormats()
    $cell.Value = ""
}

# ---------------------------------------------------------------------------
# 4. Re-shape the merged blocks for B/D/F/H/J columns around rows 28-39.
#    Previously each column merged 28:31 and 32:39 (running to 19:00); now
#    that activities stop at 17:00, they merge into a single 28:30 block.
# ---------------------------------------------------------------------------
foreach ($col in @("B", "D", "F", "H", "J")) {
    $ws.Range($col + "28:" + $col + "31").UnMerge()
    $ws.Range($col + "32:" + $col + "39").UnMerge()
    $ws.Range($col + "28:" + $col + "30").Merge()
}
